$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9.280306156861851
$ws.Range("C2").Value = 4.447645291865344
$ws.Range("D2").Value = 10.36178816753482
$ws.Range("F2").Value = 33.54911877370986
$ws.Range("G2").Value = 3.657742004949989
$ws.Range("I2").Value = 22.73886863982194
$ws.Range("J2").Value = 11.32968507126717
$ws.Range("K2").Value = 9.53859081514925
$ws.Range("M2").Value = 16.23787269867134
$ws.Range("N2").Value = 19.57683185014094
$ws.Range("O2").Value = 24.77104584977666
$ws.Range("B3").Value = 9.030210797441091
$ws.Range("C3").Value = 4.277911756317233
$ws.Range("D3").Value = 10.31951106018192
$ws.Range("F3").Value = 33.6020666184061
$ws.Range("G3").Value = 3.659404128894479
$ws.Range("I3").Value = 22.81921275519354
$ws.Range("J3").Value = 11.34716626520602
$ws.Range("K3").Value = 9.380431016574175
$ws.Range("M3").Value = 16.16713098031706
$ws.Range("N3").Value = 19.63431106048354
$ws.Range("O3").Value = 24.84136793999944
$ws.Range("B4").Value = 8.874608905566012
$ws.Range("C4").Value = 4.170760155576389
$ws.Range("D4").Value = 10.29532626544949
$ws.Range("F4").Value = 33.64169927920811
$ws.Range("G4").Value = 3.660479510867735
$ws.Range("I4").Value = 22.87230482645739
$ws.Range("J4").Value = 11.35949728460555
$ws.Range("K4").Value = 9.283472633348939
$ws.Range("M4").Value = 16.12625983289112
$ws.Range("N4").Value = 19.67126151905453
$ws.Range("O4").Value = 24.88910343816542
$ws.Range("B5").Value = 8.810784461670211
$ws.Range("C5").Value = 4.12642360227914
$ws.Range("D5").Value = 10.28592428803826
$ws.Range("F5").Value = 33.65963898400152
$ws.Range("G5").Value = 3.660931567549819
$ws.Range("I5").Value = 22.89488569068002
$ws.Range("J5").Value = 11.36492412322111
$ws.Range("K5").Value = 9.244049314856534
$ws.Range("M5").Value = 16.11026206432033
$ws.Range("N5").Value = 19.68673725808937
$ws.Range("O5").Value = 24.90970078226614
$ws.Range("B6").Value = 8.800164344842754
$ws.Range("C6").Value = 4.119023083714421
$ws.Range("D6").Value = 10.2843906989766
$ws.Range("F6").Value = 33.66272584799746
$ws.Range("G6").Value = 3.66100746774123
$ws.Range("I6").Value = 22.89869232154539
$ws.Range("J6").Value = 11.36584951962802
$ws.Range("K6").Value = 9.237509879843298
$ws.Range("M6").Value = 16.10764571703708
$ws.Range("N6").Value = 19.68933228228354
$ws.Range("O6").Value = 24.91319005665724
$ws.Range("B7").Value = 8.873749695136507
$ws.Range("C7").Value = 4.170164842694754
$ws.Range("D7").Value = 10.29519762120245
$ws.Range("F7").Value = 33.64193397937103
$ws.Range("G7").Value = 3.660485551409653
$ws.Range("I7").Value = 22.87260553218265
$ws.Range("J7").Value = 11.35956884562183
$ws.Range("K7").Value = 9.282940534336642
$ws.Range("M7").Value = 16.1260414027993
$ws.Range("N7").Value = 19.67146853569639
$ws.Range("O7").Value = 24.88937658787475
$ws.Range("B8").Value = 9.194554348399244
$ws.Range("C8").Value = 4.389767940730955
$ws.Range("D8").Value = 10.34684761622062
$ws.Range("F8").Value = 33.56589590174064
$ws.Range("G8").Value = 3.658303750435744
$ws.Range("I8").Value = 22.76579069748566
$ws.Range("J8").Value = 11.33538112581861
$ws.Range("K8").Value = 9.484055331928205
$ws.Range("M8").Value = 16.21295651526068
$ws.Range("N8").Value = 19.59630729863364
$ws.Range("O8").Value = 24.79434626784846
$ws.Range("B9").Value = 9.80342151066332
$ws.Range("C9").Value = 4.794556727772744
$ws.Range("D9").Value = 10.46185029653205
$ws.Range("F9").Value = 33.47336754297626
$ws.Range("G9").Value = 3.654458386752441
$ws.Range("I9").Value = 22.5861728044918
$ws.Range("J9").Value = 11.30061715212157
$ws.Range("K9").Value = 9.87750153812093
$ws.Range("M9").Value = 16.40315507815216
$ws.Range("N9").Value = 19.46201549633595
$ws.Range("O9").Value = 24.64420792433135
$ws.Range("B10").Value = 10.23346393059576
$ws.Range("C10").Value = 5.073221116488372
$ws.Range("D10").Value = 10.55419830277275
$ws.Range("F10").Value = 33.4399478483944
$ws.Range("G10").Value = 3.651894566194346
$ws.Range("I10").Value = 22.47241318128677
$ws.Range("J10").Value = 11.28278746354336
$ws.Range("K10").Value = 10.16319097832488
$ws.Range("M10").Value = 16.55410308086308
$ws.Range("N10").Value = 19.37125678537624
$ws.Range("O10").Value = 24.5560506286679
$ws.Range("B11").Value = 10.42442490473794
$ws.Range("C11").Value = 5.195418445343894
$ws.Range("D11").Value = 10.59779832344047
$ws.Range("F11").Value = 33.4322521946696
$ws.Range("G11").Value = 3.650784400501613
$ws.Range("I11").Value = 22.42461682156724
$ws.Range("J11").Value = 11.27634734656736
$ws.Range("K11").Value = 10.2918603607076
$ws.Range("M11").Value = 16.62501588311345
$ws.Range("N11").Value = 19.33166780094952
$ws.Range("O11").Value = 24.52076874952451
$ws.Range("B12").Value = 10.49599755445065
$ws.Range("C12").Value = 5.240999474497957
$ws.Range("D12").Value = 10.61452670961945
$ws.Range("F12").Value = 33.43041676532258
$ws.Range("G12").Value = 3.650372037776321
$ws.Range("I12").Value = 22.40708641205856
$ws.Range("J12").Value = 11.27414849623995
$ws.Range("K12").Value = 10.34035357417017
$ws.Range("M12").Value = 16.65217464703712
$ws.Range("N12").Value = 19.31691933703717
$ws.Range("O12").Value = 24.50810270161911
$ws.Range("B13").Value = 10.48061705216598
$ws.Range("C13").Value = 5.231214095476508
$ws.Range("D13").Value = 10.61091443053589
$ws.Range("F13").Value = 33.43076409900429
$ws.Range("G13").Value = 3.650460490808161
$ws.Range("I13").Value = 22.41083658288804
$ws.Range("J13").Value = 11.27461139556327
$ws.Range("K13").Value = 10.32992066937052
$ws.Range("M13").Value = 16.64631220926572
$ws.Range("N13").Value = 19.32008489123323
$ws.Range("O13").Value = 24.51079965959256
$ws.Range("B14").Value = 10.43032843554125
$ws.Range("C14").Value = 5.199182487995158
$ws.Range("D14").Value = 10.59917027120804
$ws.Range("F14").Value = 33.4320795796497
$ws.Range("G14").Value = 3.650750314398874
$ws.Range("I14").Value = 22.42316317685136
$ws.Range("J14").Value = 11.27616164051226
$ws.Range("K14").Value = 10.29585480983928
$ws.Range("M14").Value = 16.62724422665467
$ws.Range("N14").Value = 19.33044957155706
$ws.Range("O14").Value = 24.51971278252895
$ws.Range("B15").Value = 10.39942688054624
$ws.Range("C15").Value = 5.179471055119427
$ws.Range("D15").Value = 10.59200470170665
$ws.Range("F15").Value = 33.43302580111842
$ws.Range("G15").Value = 3.650928884707839
$ws.Range("I15").Value = 22.43078769259555
$ws.Range("J15").Value = 11.27714243830791
$ws.Range("K15").Value = 10.27495708446723
$ws.Range("M15").Value = 16.61560381605447
$ws.Range("N15").Value = 19.33682985446474
$ws.Range("O15").Value = 24.52526279534996
$ws.Range("B16").Value = 10.22088422205645
$ws.Range("C16").Value = 5.065140157383652
$ws.Range("D16").Value = 10.55138003612984
$ws.Range("F16").Value = 33.44060178831162
$ws.Range("G16").Value = 3.651968244307404
$ws.Range("I16").Value = 22.47561638481685
$ws.Range("J16").Value = 11.28324192941935
$ws.Range("K16").Value = 10.1547522122833
$ws.Range("M16").Value = 16.54951248425193
$ws.Range("N16").Value = 19.3738780855251
$ws.Range("O16").Value = 24.55845349906178
$ws.Range("B17").Value = 10.11010756215111
$ws.Range("C17").Value = 4.993804925079062
$ws.Range("D17").Value = 10.5268579717621
$ws.Range("F17").Value = 33.44717181330208
$ws.Range("G17").Value = 3.652620206367659
$ws.Range("I17").Value = 22.50413028761569
$ws.Range("J17").Value = 11.28741143622078
$ws.Range("K17").Value = 10.08064659622925
$ws.Range("M17").Value = 16.50953062930657
$ws.Range("N17").Value = 19.39703999553204
$ws.Range("O17").Value = 24.58005061627078
$ws.Range("B18").Value = 10.04595614681916
$ws.Range("C18").Value = 4.952346571375758
$ws.Range("D18").Value = 10.51290389353073
$ws.Range("F18").Value = 33.45165740898122
$ws.Range("G18").Value = 3.653000483114572
$ws.Range("I18").Value = 22.52090281664815
$ws.Range("J18").Value = 11.28996690555029
$ws.Range("K18").Value = 10.0379039626787
$ws.Range("M18").Value = 16.48674657814317
$ws.Range("N18").Value = 19.41052196489139
$ws.Range("O18").Value = 24.5929264506913
$ws.Range("B19").Value = 10.02416296787823
$ws.Range("C19").Value = 4.938237101078827
$ws.Range("D19").Value = 10.50820542790454
$ws.Range("F19").Value = 33.45329754298977
$ws.Range("G19").Value = 3.65313014723655
$ws.Range("I19").Value = 22.52664559693589
$ws.Range("J19").Value = 11.29085916651008
$ws.Range("K19").Value = 10.02341300879578
$ws.Range("M19").Value = 16.47906929886139
$ws.Range("N19").Value = 19.41511421604094
$ws.Range("O19").Value = 24.59736388165047
$ws.Range("B20").Value = 10.12194553021787
$ws.Range("C20").Value = 5.001443273547243
$ws.Range("D20").Value = 10.52945290160826
$ws.Range("F20").Value = 33.44639928668321
$ws.Range("G20").Value = 3.652550257152012
$ws.Range("I20").Value = 22.50105642061964
$ws.Range("J20").Value = 11.28695130972619
$ws.Range("K20").Value = 10.08854792062184
$ws.Range("M20").Value = 16.51376490172343
$ws.Range("N20").Value = 19.39455783457115
$ws.Range("O20").Value = 24.57770459692729
$ws.Range("B21").Value = 10.4451200264158
$ws.Range("C21").Value = 5.208609999504078
$ws.Range("D21").Value = 10.60261398587513
$ws.Range("F21").Value = 33.43166392334074
$ws.Range("G21").Value = 3.650664968477086
$ws.Range("I21").Value = 22.41952710969122
$ws.Range("J21").Value = 11.27569978883371
$ws.Range("K21").Value = 10.30586740070635
$ws.Range("M21").Value = 16.63283680113659
$ws.Range("N21").Value = 19.32739862669017
$ws.Range("O21").Value = 24.51707592701344
$ws.Range("B22").Value = 10.65199162204061
$ws.Range("C22").Value = 5.339956141432549
$ws.Range("D22").Value = 10.65169448547679
$ws.Range("F22").Value = 33.42832049030286
$ws.Range("G22").Value = 3.649479627827755
$ws.Range("I22").Value = 22.36956001345859
$ws.Range("J22").Value = 11.26974427713712
$ws.Range("K22").Value = 10.44653178025833
$ws.Range("M22").Value = 16.71243023593575
$ws.Range("N22").Value = 19.2849222416709
$ws.Range("O22").Value = 24.48149972933265
$ws.Range("B23").Value = 10.54199878152462
$ws.Range("C23").Value = 5.27023525268085
$ws.Range("D23").Value = 10.62538712699751
$ws.Range("F23").Value = 33.42953009624461
$ws.Range("G23").Value = 3.650107996420468
$ws.Range("I23").Value = 22.39592470769854
$ws.Range("J23").Value = 11.27279506277905
$ws.Range("K23").Value = 10.37159599658481
$ws.Range("M23").Value = 16.66979332544007
$ws.Range("N23").Value = 19.30746347703055
$ws.Range("O23").Value = 24.5001166631949
$ws.Range("B24").Value = 10.11659502616495
$ws.Range("C24").Value = 4.997991365073335
$ws.Range("D24").Value = 10.52827928453103
$ws.Range("F24").Value = 33.44674633929606
$ws.Range("G24").Value = 3.652581864197055
$ws.Range("I24").Value = 22.50244493378303
$ws.Range("J24").Value = 11.28715883959905
$ws.Range("K24").Value = 10.08497615924971
$ws.Range("M24").Value = 16.5118499564888
$ws.Range("N24").Value = 19.39567950274713
$ws.Range("O24").Value = 24.57876380146191
$ws.Range("B25").Value = 9.641409246736805
$ws.Range("C25").Value = 4.688149742261509
$ws.Range("D25").Value = 10.42932097683448
$ws.Range("F25").Value = 33.49233231849311
$ws.Range("G25").Value = 3.655452567668083
$ws.Range("I25").Value = 22.63156793804531
$ws.Range("J25").Value = 11.30866649331136
$ws.Range("K25").Value = 9.771459093436958
$ws.Range("M25").Value = 16.34966971815853
$ws.Range("N25").Value = 19.49695095320811
$ws.Range("O25").Value = 24.68093947236139
